$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Förändrad" (column C) dates for rows 2-5 from 2023-09-01 (45170) to 2023-09-05 (45174)
$ws.Range("C2").Value = 45174
$ws.Range("C3").Value = 45174
$ws.Range("C4").Value = 45174
$ws.Range("C5").Value = 45174
